$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'255.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'4.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'27.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'-8.06%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'5.190"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'-0.42%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'0.05856"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'1.92%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'6.721"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'1.14%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.8694"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'1.32%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'0.9633"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'12.98%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.1408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'2.02%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Formula = "'1.17%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.03203"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'1.85%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'0.09223"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'-1.40%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'0.001546"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'1.36%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Formula = "'-94.03%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.005999"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'-0.26%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Formula = "'-0.77%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Formula = "'-1.29%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Formula = "'1.50%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'0.68%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'0.03449"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'3.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'0.1281"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'-2.08%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'3.522"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'1.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.04184"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'1.67%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Formula = "'-0.77%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Formula = "'0.004790"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'15.00%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'-0.01%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Formula = "'0.03817"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'1.48%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.005638"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'-0.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Formula = "'2.98%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'0.002299"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'4.54%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'0.009832"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'-4.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.00005423"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'8.26%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'-0.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Formula = "'11.26%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'0.002128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'-3.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'-0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'-0.01%"
$ws.Range("E50").Style = "Normal"
